$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.950.04"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "2.249.99"
$ws.Range("E3").Value = "  +2.59%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.85"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.95"
$ws.Range("E6").Value = "  +1.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0948"
$ws.Range("E9").Value = "  +2.94%  "

$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("E11").Value = "  +3.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.64"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").Value = "2.649.07"
$ws.Range("E13").Value = "  +2.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.68"
$ws.Range("E14").Value = "  +2.88%  "

$ws.Range("D15").Value = "53.857.24"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").Value = "2.229.25"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.99"
$ws.Range("E18").Value = "  +4.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.06"
$ws.Range("E19").Value = "  +2.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "299.69"
$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.40"
$ws.Range("E21").Value = "  +4.66%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -2.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.86"
$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.369"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("D27").Value = "2.353.90"
$ws.Range("E27").Value = "  +2.61%  "

$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.07"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.60"
$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0675"
$ws.Range("E32").Value = "  +1.88%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  +2.40%  "

$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.61"
$ws.Range("E37").Value = "  +2.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  +7.11%  "

$ws.Range("E39").Value = "  +2.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  +3.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.67"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.38"
$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.370"
$ws.Range("E43").Value = "  +1.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  +2.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.91"
$ws.Range("E45").Value = "  +3.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.65"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0884"
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.537"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "234.39"
$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0480"
$ws.Range("E50").Value = "  +2.02%  "

$ws.Range("E51").Value = "  +0.45%  "
